$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.597.79"
$ws.Range("E2").Value = "  +3.71%  "

$ws.Range("D3").Value = "3.255.08"
$ws.Range("E3").Value = "  +6.83%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.43"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.77"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "3.247.96"
$ws.Range("E8").Value = "  +6.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.05"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.78%  "

$ws.Range("E11").Value = "  +5.74%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.491"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.00"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "

$ws.Range("E14").Value = "  +5.58%  "

$ws.Range("D15").Value = "3.780.84"
$ws.Range("E15").Value = "  +7.23%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.650.65"
$ws.Range("E16").Value = "  +3.68%  "

$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "556.98"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.54%  "

$ws.Range("D18").Value = "3.256.95"
$ws.Range("E18").Value = "  +7.07%  "

$ws.Range("E19").Value = "  +3.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.70%  "

$ws.Range("E21").Value = "  +4.44%  "

$ws.Range("E22").Value = "  +7.99%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.83"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.92%  "

$ws.Range("E24").Value = "  +7.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.96"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.95%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.32"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +18.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.98"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "27.87"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.76"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.90%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "564.10"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.81%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.71"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.42"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0460"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.40"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.70%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0866"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.22%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.04"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.51%  "

$ws.Range("D42").Value = "3.196.66"
$ws.Range("E42").Value = "  +9.78%  "

$ws.Range("E43").Value = "  +2.61%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.42%  "

$ws.Range("E45").Value = "  +9.32%  "

$ws.Range("E46").Value = "  +4.28%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "

$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.0₃0560"
$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.53"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.21%  "

$ws.Range("E50").Value = "  +2.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.24"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.90%  "
